# Saldo.xlsx update:
#  - Remove a batch of account rows that dropped out of the export.
#  - Re-add account 004313254 (GUSTAVO) with an updated balance (39137.61)
#    in a new position (just above account 004254210 / MARCO).
#
# Row numbers below are the ORIGINAL (before-edit) 1-based worksheet rows
# that must be deleted; they are processed highest-to-lowest so earlier
# deletions never shift the row number of one still to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(20, 19, 18, 17, 16, 15, 12, 11, 10, 9, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# After the deletions above, account 001882235 (LAGO) sits on row 4 and
# account 004254210 (MARCO) sits on row 5. Insert a fresh row 5 for the
# re-added GUSTAVO record, ahead of MARCO.
$ws.Rows(5).Insert()

# Leading apostrophe forces text storage so the leading zeros in the
# account number survive (matches the other inline-string account cells).
$ws.Range("A5").Value = "'004313254"
$ws.Range("B5").Value = "GUSTAVO"
$ws.Range("C5").Value = 39137.61
